$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("A2").Value = 71362
$ws.Range("A3").Value = 71405
$ws.Range("A4").Value = 71409

# Remove rows 5 through 8 (old extra data), shrinking the used range to A1:A4
$ws.Range("A5:A8").ClearContents()

# Update the selected cell/range to match the new state
$ws.Range("E9").Select()
